$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 7450
$ws.Range("J7").Value = 7450
$ws.Range("L7").Value = 7450
$ws.Range("N7").Value = -7674
$ws.Range("H14").Value = 7450
$ws.Range("J14").Value = 7450
$ws.Range("L14").Value = 7450
$ws.Range("N14").Value = -7832
$ws.Range("H53").Value = 961
$ws.Range("I53").Value = 1002.875
$ws.Range("K53").Value = 1002.875
$ws.Range("M53").Value = -365.875
$ws.Range("H88").Value = 5265892.5
$ws.Range("I88").Value = 11113120
$ws.Range("J88").Value = 3387.2
$ws.Range("K88").Value = 11113120
$ws.Range("L88").Value = 3387.2
$ws.Range("M88").Value = -11112714
$ws.Range("N88").Value = -4199.2
$ws.Range("H91").Value = 5265892.5
$ws.Range("I91").Value = 11113120
$ws.Range("J91").Value = 3387.2
$ws.Range("K91").Value = 11113120
$ws.Range("L91").Value = 3387.2
$ws.Range("M91").Value = -11111716
$ws.Range("N91").Value = -6195.2
$ws.Range("H107").Value = 567.35297
$ws.Range("I107").Value = 477.8125
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 477.8125
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1442.1875
$ws.Range("N107").Value = -5840
$ws.Range("H129").Value = 1316.6957
$ws.Range("I129").Value = 479
$ws.Range("J129").Value = 2887.375
$ws.Range("K129").Value = 1437
$ws.Range("L129").Value = 8662.125
$ws.Range("M129").Value = 3563
$ws.Range("N129").Value = -18662.125
$ws.Range("H132").Value = 5915.864
$ws.Range("I132").Value = 6126.1904
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 18378.5712
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -15848.5712
$ws.Range("N132").Value = -9557
$ws.Range("H137").Value = 1614548.9
$ws.Range("I137").Value = 2274492.8
$ws.Range("J137").Value = 1352.6666
$ws.Range("K137").Value = 6823478.399999999
$ws.Range("L137").Value = 4057.9998
$ws.Range("M137").Value = -6820928.399999999
$ws.Range("N137").Value = -9157.9998
$ws.Range("H141").Value = 2319.5334
$ws.Range("I141").Value = 2319.5334
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6958.600199999999
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -1778.600199999999

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2732267
$ws.Range("I32").Value = 1373268.5
$ws.Range("K32").Value = 1373268.5
$ws.Range("M32").Value = -1372981.5
$ws.Range("H74").Value = 1601.7667
$ws.Range("I74").Value = 1059.0952
$ws.Range("K74").Value = 1059.0952
$ws.Range("M74").Value = -185.0952
$ws.Range("H77").Value = 1601.7667
$ws.Range("I77").Value = 1059.0952
$ws.Range("K77").Value = 5295.476
$ws.Range("M77").Value = -927.4759999999997
$ws.Range("H80").Value = 89000
$ws.Range("I80").Value = 78000
$ws.Range("K80").Value = 78000
$ws.Range("M80").Value = -77002
$ws.Range("H83").Value = 89000
$ws.Range("I83").Value = 78000
$ws.Range("K83").Value = 234000
$ws.Range("M83").Value = -229008
$ws.Range("H97").Value = 650.2353000000001
$ws.Range("I97").Value = 546.7692
$ws.Range("K97").Value = 546.7692
$ws.Range("M97").Value = -50.76919999999996

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3561.639
$ws.Range("I86").Value = 3565.08
$ws.Range("J86").Value = 3553.818
$ws.Range("K86").Value = 3565.08
$ws.Range("L86").Value = 3553.818
$ws.Range("M86").Value = -2442.08
$ws.Range("N86").Value = -5799.818
$ws.Range("H89").Value = 3561.639
$ws.Range("I89").Value = 3565.08
$ws.Range("J89").Value = 3553.818
$ws.Range("K89").Value = 17825.4
$ws.Range("L89").Value = 17769.09
$ws.Range("M89").Value = -12209.4
$ws.Range("N89").Value = -29001.09

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1123.909
$ws.Range("I16").Value = 1040.5555
$ws.Range("K16").Value = 1040.5555
$ws.Range("M16").Value = -753.5554999999999
$ws.Range("H58").Value = 1915.9445
$ws.Range("I58").Value = 1435.2222
$ws.Range("K58").Value = 1435.2222
$ws.Range("M58").Value = -1232.2222
$ws.Range("H62").Value = 499.5
$ws.Range("I62").Value = 500
$ws.Range("J62").Value = 499
$ws.Range("K62").Value = 500
$ws.Range("L62").Value = 499
$ws.Range("M62").Value = 124
$ws.Range("N62").Value = -1747
$ws.Range("H65").Value = 499.5
$ws.Range("I65").Value = 500
$ws.Range("J65").Value = 499
$ws.Range("K65").Value = 2500
$ws.Range("L65").Value = 2495
$ws.Range("M65").Value = 620
$ws.Range("N65").Value = -8735
$ws.Range("H113").Value = 1123.909
$ws.Range("I113").Value = 1040.5555
$ws.Range("K113").Value = 1040.5555
$ws.Range("M113").Value = 1129.4445
$ws.Range("H132").Value = 3647.65
$ws.Range("I132").Value = 2969.7666
$ws.Range("K132").Value = 8909.299800000001
$ws.Range("M132").Value = -6379.299800000001
$ws.Range("H134").Value = 3734.973
$ws.Range("I134").Value = 3939.9666
$ws.Range("K134").Value = 11819.8998
$ws.Range("M134").Value = -9284.899800000001
$ws.Range("H136").Value = 1915.9445
$ws.Range("I136").Value = 1435.2222
$ws.Range("K136").Value = 4305.6666
$ws.Range("M136").Value = -1755.6666

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1477195.6
$ws.Range("I131").Value = 25297.375
$ws.Range("J131").Value = 1840170.2
$ws.Range("K131").Value = 75892.125
$ws.Range("L131").Value = 5520510.6
$ws.Range("M131").Value = -70852.125
$ws.Range("N131").Value = -5530590.6

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 776.08
$ws.Range("I97").Value = 826.7273
$ws.Range("K97").Value = 826.7273
$ws.Range("M97").Value = -330.7273
$ws.Range("H132").Value = 2394.3572
$ws.Range("J132").Value = 2544.2222
$ws.Range("L132").Value = 7632.6666
$ws.Range("N132").Value = -12692.6666

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H100").Value = 2436.125
$ws.Range("I100").Value = 2938
$ws.Range("J100").Value = 1599.6666
$ws.Range("K100").Value = 2938
$ws.Range("L100").Value = 1599.6666
$ws.Range("M100").Value = -2397
$ws.Range("N100").Value = -2681.6666
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 5082.4287
$ws.Range("I132").Value = 5633.25
$ws.Range("K132").Value = 16899.75
$ws.Range("M132").Value = -14369.75
$ws.Range("H136").Value = 4879.591
$ws.Range("I136").Value = 4908.1113
$ws.Range("K136").Value = 14724.3339
$ws.Range("M136").Value = -12174.3339

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 780
$ws.Range("I2").Value = 755.55554
$ws.Range("K2").Value = 755.55554
$ws.Range("M2").Value = -643.55554
$ws.Range("H4").Value = 296.73334
$ws.Range("I4").Value = 191.58333
$ws.Range("J4").Value = 717.3333
$ws.Range("K4").Value = 191.58333
$ws.Range("L4").Value = 717.3333
$ws.Range("M4").Value = -78.58332999999999
$ws.Range("N4").Value = -943.3333
$ws.Range("H14").Value = 10111
$ws.Range("I14").Value = 10111
$ws.Range("K14").Value = 10111
$ws.Range("M14").Value = -9943
$ws.Range("H22").Value = 8503
$ws.Range("J22").Value = 9666.333000000001
$ws.Range("L22").Value = 9666.333000000001
$ws.Range("N22").Value = -10252.333
$ws.Range("H132").Value = 3221.0952
$ws.Range("I132").Value = 3309.875
$ws.Range("J132").Value = 2937
$ws.Range("K132").Value = 9929.625
$ws.Range("L132").Value = 8811
$ws.Range("M132").Value = -7399.625
$ws.Range("N132").Value = -13871
$ws.Range("H136").Value = 2844.9375
$ws.Range("I136").Value = 3192.4285
$ws.Range("K136").Value = 9577.2855
$ws.Range("M136").Value = -7027.2855

Write-Host "All edits applied."